# Add a new list item ("آب") at the very end of the document, after the
# last existing "ListParagraph" (numId 2) item, matching the list's
# formatting (sz 22 / szCs 24 / rtl / cs-hinted font).

$d = $word.ActiveDocument

# Move to the end of the document and insert a brand-new paragraph there.
# Word will carry the paragraph/run formatting of the preceding list item
# forward onto the new (empty) paragraph.
$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter()

# Collapse to the (now) very end of the document -- i.e. inside the newly
# created empty paragraph -- and give it its text.
$end = $d.Content
$end.Collapse(0)
$end.Text = "آب"
